$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.030.43"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.238.24"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'305.15"
$ws.Range("E5").Value = "  -4.77%  "
$ws.Range("D6").Value = "'95.79"
$ws.Range("E6").Value = "  -6.56%  "
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  -5.76%  "
$ws.Range("D10").Value = "'34.99"
$ws.Range("E10").Value = "  -6.51%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("D12").Value = "'7.22"
$ws.Range("E12").Value = "  -5.74%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "2.579.11"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "2.239.17"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "'0.827"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D17").Value = "'13.56"
$ws.Range("E17").Value = "  -6.56%  "
$ws.Range("D18").Value = "43.905.91"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "0.0₃0958"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "'12.29"
$ws.Range("E20").Value = "  -8.92%  "
$ws.Range("D21").Value = "'6.23"
$ws.Range("E21").Value = "  -5.17%  "
$ws.Range("D22").Value = "'64.77"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "'236.32"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'2.92"
$ws.Range("E24").Value = "  -7.66%  "
$ws.Range("E25").Value = "  -7.80%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'9.93"
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").Value = "'37.73"
$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "'5.96"
$ws.Range("E30").Value = "  -5.05%  "
$ws.Range("D31").Value = "'20.04"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").Value = "'155.59"
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("D33").Value = "'0.0808"
$ws.Range("E33").Value = "  -5.48%  "
$ws.Range("E34").Value = "  +6.24%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  -11.24%  "
$ws.Range("D39").Value = "'15.35"
$ws.Range("E39").Value = "  -8.11%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'3.82"
$ws.Range("E40").Value = "  -9.86%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "'3.35"
$ws.Range("E41").Value = "  -10.06%  "
$ws.Range("E42").Value = "  -5.43%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "1.738.92"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "'85.34"
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").Value = "'0.189"
$ws.Range("E46").Value = "  -5.54%  "
$ws.Range("D47").Value = "'100.05"
$ws.Range("E47").Value = "  -5.03%  "
$ws.Range("D48").Value = "'4.94"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("D49").Value = "'69.33"
$ws.Range("E49").Value = "  -7.81%  "
$ws.Range("D50").Value = "'8.09"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").Value = "'54.40"
$ws.Range("E51").Value = "  -7.22%  "
